# Generate Report for Handback
# Replaces the two handback-item identifiers (GUID-ish filenames) and their
# associated timestamps across the Overview / zh-cn / de-de sheets, matching
# a freshly regenerated handback-status report.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$oldId1 = "331089a6-3e31-40ba-9acf-340eb7ecd295"
$newId1 = "fa547c0d-cd1c-4e34-8711-6f34bda7bb82"
$oldId2 = "b5896939-8c8b-4b80-88b2-65bab88e3381"
$newId2 = "ffffc7cc1db2-99b0-4511-bbb1-1592252a3feb"

# ---------------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------------
$wsOverview.Range("A2").Value = "$newId1.md"
$wsOverview.Range("G2").Value = "2016-08-25 19:05:18"

$wsOverview.Range("A3").Value = "$newId2.md"
$wsOverview.Range("G3").Value = "2016-08-25 19:05:18"

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$wsZhCn.Range("A2").Value = "$newId1.md"
$wsZhCn.Range("G2").Value = "$newId1.dec929f9212e255ca06cd666eafdecb727f7f07f.zh-cn.xlf"
$wsZhCn.Range("H2").Value = "2016-08-25 19:05:13"
$wsZhCn.Range("I2").Value = "$newId1.md"
$wsZhCn.Range("J2").Value = "$newId1.dec929f9212e255ca06cd666eafdecb727f7f07f.zh-cn.xlf"
$wsZhCn.Range("K2").Value = "2016-08-25 19:05:49"

$wsZhCn.Range("A3").Value = "$newId2.md"
$wsZhCn.Range("G3").Value = "$newId1.dec929f9212e255ca06cd666eafdecb727f7f07f.zh-cn.xlf"
$wsZhCn.Range("H3").Value = "2016-08-25 19:05:13"
$wsZhCn.Range("I3").Value = "$newId2.md"
$wsZhCn.Range("J3").Value = "$newId1.dec929f9212e255ca06cd666eafdecb727f7f07f.zh-cn.xlf"
$wsZhCn.Range("K3").Value = "2016-08-25 19:05:49"

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$wsDeDe.Range("A2").Value = "$newId1.md"
$wsDeDe.Range("G2").Value = "$newId1.dec929f9212e255ca06cd666eafdecb727f7f07f.de-de.xlf"
$wsDeDe.Range("H2").Value = "2016-08-25 19:05:18"
$wsDeDe.Range("I2").Value = "$newId1.md"
$wsDeDe.Range("J2").Value = "$newId1.dec929f9212e255ca06cd666eafdecb727f7f07f.de-de.xlf"
$wsDeDe.Range("K2").Value = "2016-08-25 19:05:55"

$wsDeDe.Range("A3").Value = "$newId2.md"
$wsDeDe.Range("G3").Value = "$newId1.dec929f9212e255ca06cd666eafdecb727f7f07f.de-de.xlf"
$wsDeDe.Range("H3").Value = "2016-08-25 19:05:18"
$wsDeDe.Range("I3").Value = "$newId2.md"
$wsDeDe.Range("J3").Value = "$newId1.dec929f9212e255ca06cd666eafdecb727f7f07f.de-de.xlf"
$wsDeDe.Range("K3").Value = "2016-08-25 19:05:55"

# ---------------------------------------------------------------------------
# Hyperlinks: keep the original link targets (they are not touched by this
# commit), only refresh the cell text / display strings to the new names.
# Hyperlinks.Add() re-creates the link in place of the deleted one with the
# same address, so relationship ids are regenerated in the same order
# (rId2, rId3, ...) as before.
# ---------------------------------------------------------------------------

$overviewUrl1 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/cd9ef086b4143d77286273a7dca6c06f8de709e4/e2e/$oldId1.md"
$overviewUrl2 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/cd9ef086b4143d77286273a7dca6c06f8de709e4/e2e/$oldId2.md"

$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), $overviewUrl1, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "e2e\$newId1.md")
$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), $overviewUrl2, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "e2e\$newId2.md")

$zhcnUrlA1 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/cd9ef086b4143d77286273a7dca6c06f8de709e4/e2e/$oldId1.md"
$zhcnUrlI1 = "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/5edb393fc8ebf8616c9391f7d75c937df37a054c/e2e/$oldId1.md"
$zhcnUrlA2 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/cd9ef086b4143d77286273a7dca6c06f8de709e4/e2e/$oldId2.md"
$zhcnUrlI2 = "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/5edb393fc8ebf8616c9391f7d75c937df37a054c/e2e/$oldId2.md"

$wsZhCn.Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), $zhcnUrlA1, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "$newId1.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I2"), $zhcnUrlI1, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "$newId1.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), $zhcnUrlA2, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "$newId2.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I3"), $zhcnUrlI2, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "$newId2.md")

$dedeUrlA1 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/cd9ef086b4143d77286273a7dca6c06f8de709e4/e2e/$oldId1.md"
$dedeUrlI1 = "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/71aa6c7e65dd2f1bde82d8c28085bbb32567ed7f/e2e/$oldId1.md"
$dedeUrlA2 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/cd9ef086b4143d77286273a7dca6c06f8de709e4/e2e/$oldId2.md"
$dedeUrlI2 = "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/71aa6c7e65dd2f1bde82d8c28085bbb32567ed7f/e2e/$oldId2.md"

$wsDeDe.Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), $dedeUrlA1, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "$newId1.md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I2"), $dedeUrlI1, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "$newId1.md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), $dedeUrlA2, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "$newId2.md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I3"), $dedeUrlI2, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "$newId2.md")

# ---------------------------------------------------------------------------
# Hyperlinks.Add() re-applies Excel's built-in "Hyperlink" look (underlined,
# themed blue) to the touched cells. Restore the workbook's own custom
# "HyperLink" text formatting (underlined, RGB 0x6495ED) so the cells keep
# looking the way they did before the link was refreshed.
# ---------------------------------------------------------------------------
$linkCells = @(
    $wsOverview.Range("B2"), $wsOverview.Range("B3"),
    $wsZhCn.Range("A2"), $wsZhCn.Range("I2"), $wsZhCn.Range("A3"), $wsZhCn.Range("I3"),
    $wsDeDe.Range("A2"), $wsDeDe.Range("I2"), $wsDeDe.Range("A3"), $wsDeDe.Range("I3")
)
foreach ($cell in $linkCells) {
    $cell.Font.Name = "Calibri"
    $cell.Font.Size = 11
    $cell.Font.Underline = $true
    $cell.Font.Color = 15570276
}
